$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment") entirely; this shifts C->B, D->C, E->D, F->E
$ws.Columns.Item(2).Delete()

# Update header row texts for the remaining columns (now B, C, D, E)
$ws.Range("B1").Value = "All.global"
$ws.Range("C1").Value = "Males.global"
$ws.Range("D1").Value = "Females.global"
$ws.Range("E1").Value = "Not known / missing.global"
